$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 355
$ws.Range("I2").Value = 349.375
$ws.Range("K2").Value = 349.375
$ws.Range("M2").Value = -236.375
# Row 33
$ws.Range("H33").Value = 114.3125
$ws.Range("I33").Value = 102.57143
$ws.Range("J33").Value = 196.5
$ws.Range("K33").Value = 102.57143
$ws.Range("L33").Value = 196.5
$ws.Range("M33").Value = 126.42857
$ws.Range("N33").Value = -654.5
# Row 121
$ws.Range("H121").Value = 1565.4166
$ws.Range("J121").Value = 1735
$ws.Range("L121").Value = 5205
$ws.Range("N121").Value = -8699
# Row 137
$ws.Range("H137").Value = 1480.317
$ws.Range("I137").Value = 1416.56
$ws.Range("K137").Value = 4249.68
$ws.Range("M137").Value = -1699.68
# Row 138
$ws.Range("H138").Value = 2534.9773
$ws.Range("I138").Value = 1583.3729
$ws.Range("J138").Value = 4471
$ws.Range("K138").Value = 4750.1187
$ws.Range("L138").Value = 13413
$ws.Range("M138").Value = 389.8813
$ws.Range("N138").Value = -23693
# Row 141
$ws.Range("H141").Value = 4696.8975
$ws.Range("I141").Value = 2110.389
$ws.Range("J141").Value = 35735
$ws.Range("K141").Value = 6331.167
$ws.Range("L141").Value = 107205
$ws.Range("M141").Value = -1151.167
$ws.Range("N141").Value = -117565

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 882.13043
$ws.Range("I61").Value = 810.6486
$ws.Range("K61").Value = 810.6486
$ws.Range("M61").Value = -598.6486
# Row 74
$ws.Range("H74").Value = 1050.76
$ws.Range("I74").Value = 1003.2857
$ws.Range("K74").Value = 1003.2857
$ws.Range("M74").Value = -129.2857
# Row 77
$ws.Range("H77").Value = 1050.76
$ws.Range("I77").Value = 1003.2857
$ws.Range("K77").Value = 5016.4285
$ws.Range("M77").Value = -648.4285
# Row 102
$ws.Range("H102").Value = 72759.21000000001
$ws.Range("I102").Value = 1433
$ws.Range("J102").Value = 1000000
$ws.Range("K102").Value = 1433
$ws.Range("L102").Value = 1000000
$ws.Range("M102").Value = 189
$ws.Range("N102").Value = -1003244
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("N104").ClearContents()
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").ClearContents()
# Row 132
$ws.Range("H132").Value = 1607.9783
$ws.Range("I132").Value = 1116.7142
$ws.Range("J132").Value = 2372.1667
$ws.Range("K132").Value = 3350.1426
$ws.Range("L132").Value = 7116.500100000001
$ws.Range("M132").Value = -820.1425999999997
$ws.Range("N132").Value = -12176.5001
# Row 136
$ws.Range("H136").Value = 882.13043
$ws.Range("I136").Value = 810.6486
$ws.Range("K136").Value = 2431.9458
$ws.Range("M136").Value = 118.0542

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1688.375
$ws.Range("I99").Value = 921.5
$ws.Range("K99").Value = 921.5
$ws.Range("M99").Value = 576.5
# Row 134
$ws.Range("H134").Value = 2846.606
$ws.Range("I134").Value = 2288.652
$ws.Range("J134").Value = 4129.9
$ws.Range("K134").Value = 6865.956
$ws.Range("L134").Value = 12389.7
$ws.Range("M134").Value = -4330.956
$ws.Range("N134").Value = -17459.7

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1030311.9
$ws.Range("I58").Value = 1950531.8
$ws.Range("J58").Value = 1830.8235
$ws.Range("K58").Value = 1950531.8
$ws.Range("L58").Value = 1830.8235
$ws.Range("M58").Value = -1950328.8
$ws.Range("N58").Value = -2236.8235
# Row 62
$ws.Range("H62").Value = 73457.86
$ws.Range("I62").Value = 101931
$ws.Range("J62").Value = 2275
$ws.Range("K62").Value = 101931
$ws.Range("L62").Value = 2275
$ws.Range("M62").Value = -101307
$ws.Range("N62").Value = -3523
# Row 65
$ws.Range("H65").Value = 73457.86
$ws.Range("I65").Value = 101931
$ws.Range("J65").Value = 2275
$ws.Range("K65").Value = 509655
$ws.Range("L65").Value = 11375
$ws.Range("M65").Value = -506535
$ws.Range("N65").Value = -17615
# Row 100
$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164
# Row 130
$ws.Range("H130").Value = 25000
$ws.Range("J130").Value = 25000
$ws.Range("L130").Value = 25000
$ws.Range("N130").Value = -35040
# Row 132
$ws.Range("H132").Value = 376763.6
$ws.Range("I132").Value = 501623.53
$ws.Range("J132").Value = 2183.7778
$ws.Range("K132").Value = 1504870.59
$ws.Range("L132").Value = 6551.3334
$ws.Range("M132").Value = -1502340.59
$ws.Range("N132").Value = -11611.3334
# Row 134
$ws.Range("H134").Value = 1386.5883
$ws.Range("I134").Value = 1091.0834
$ws.Range("J134").Value = 2095.8
$ws.Range("K134").Value = 3273.2502
$ws.Range("L134").Value = 6287.400000000001
$ws.Range("M134").Value = -738.2501999999999
$ws.Range("N134").Value = -11357.4
# Row 136
$ws.Range("H136").Value = 1030311.9
$ws.Range("I136").Value = 1950531.8
$ws.Range("J136").Value = 1830.8235
$ws.Range("K136").Value = 5851595.4
$ws.Range("L136").Value = 5492.470499999999
$ws.Range("M136").Value = -5849045.4
$ws.Range("N136").Value = -10592.4705

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 528.625
$ws.Range("I113").Value = 440.2
$ws.Range("J113").Value = 676
$ws.Range("K113").Value = 1320.6
$ws.Range("L113").Value = 2028
$ws.Range("M113").Value = 849.4000000000001
$ws.Range("N113").Value = -6368

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1721.025
$ws.Range("I132").Value = 1117.7667
$ws.Range("J132").Value = 3530.8
$ws.Range("K132").Value = 3353.300099999999
$ws.Range("L132").Value = 10592.4
$ws.Range("M132").Value = -823.3000999999995
$ws.Range("N132").Value = -15652.4

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 4000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
# Row 129
$ws.Range("H129").Value = 29000
$ws.Range("J129").Value = 29000
$ws.Range("L129").Value = 29000
$ws.Range("N129").Value = -39000
# Row 136
$ws.Range("H136").Value = 3000.2354
$ws.Range("I136").Value = 2540.12
$ws.Range("K136").Value = 7620.36
$ws.Range("M136").Value = -5070.36

$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 49205.2
$ws.Range("I26").Value = 44003
$ws.Range("J26").Value = 70014
$ws.Range("K26").Value = 44003
$ws.Range("L26").Value = 70014
$ws.Range("M26").Value = -43710
$ws.Range("N26").Value = -70600
# Row 107
$ws.Range("H107").Value = 545.64
$ws.Range("I107").Value = 501.9
$ws.Range("K107").Value = 1505.7
$ws.Range("M107").Value = 414.3000000000002
# Row 132
$ws.Range("H132").Value = 921.3269
$ws.Range("I132").Value = 787
$ws.Range("J132").Value = 1485.5
$ws.Range("K132").Value = 2361
$ws.Range("L132").Value = 4456.5
$ws.Range("M132").Value = 169
$ws.Range("N132").Value = -9516.5
